# The source export re-ordered two pairs of observation records:
#   - row 5  <->  row 6   (a "Lunglav"/"Granticka" pair)
#   - row 7  <->  row 8   (two "Tretaig hackspett" records)
# Every field of one row trades places with the matching field of its
# partner row. Only the columns whose content actually differs between
# the two rows of a pair are touched, so columns that are already
# identical in both rows (dates, location, observer, ...) are left
# completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Cell($ws, $row1, $row2, $col) {
    $cell1 = $ws.Cells.Item($row1, $col)
    $cell2 = $ws.Cells.Item($row2, $col)

    $val1 = $cell1.Value2
    $val2 = $cell2.Value2

    # Column I ("Antal") stores small counts as plain digit text in the
    # source file (e.g. "1"); keep that as text rather than letting
    # Excel reinterpret it as a number when it lands in the other row.
    if ($col -eq 9) {
        if ($val2 -eq $null -or $val2 -eq "") {
            $cell1.Value2 = $val2
        } else {
            $cell1.Value2 = "'" + $val2
        }
        if ($val1 -eq $null -or $val1 -eq "") {
            $cell2.Value2 = $val1
        } else {
            $cell2.Value2 = "'" + $val1
        }
    } else {
        $cell1.Value2 = $val2
        $cell2.Value2 = $val1
    }
}

# Columns A,B,E,F,G,H,Q,R differ between row 5 and row 6 -> swap them.
$cols56 = @(1, 2, 5, 6, 7, 8, 17, 18)
foreach ($c in $cols56) {
    Swap-Cell $ws 5 6 $c
}

# Columns A,I,M,N,Q,R,AC,AE differ between row 7 and row 8 -> swap them.
$cols78 = @(1, 9, 13, 14, 17, 18, 29, 31)
foreach ($c in $cols78) {
    Swap-Cell $ws 7 8 $c
}
